$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disambiguate the mXXX_xxx #defines used as "row" markers (C/I/J columns)
# by renaming them to mROW_* variants, distinct from the mPOWERON/mMENU/
# mOPNBRL/mLOKLOD #defines used elsewhere (columns D/E "effect" cells keep
# their original mEFCT_* names and are untouched).

$ws.Range("C2").Value = "mROW_POWERON"

$ws.Range("J2").Value = "mROW_MENU"
$ws.Range("C4").Value = "mROW_MENU"
$ws.Range("I4").Value = "mROW_MENU"
$ws.Range("C5").Value = "mROW_MENU"
$ws.Range("C6").Value = "mROW_MENU"
$ws.Range("J8").Value = "mROW_MENU"
$ws.Range("J10").Value = "mROW_MENU"

$ws.Range("I5").Value = "mROW_OPNBRL"
$ws.Range("C8").Value = "mROW_OPNBRL"

$ws.Range("I6").Value = "mROW_LOKLOD"
$ws.Range("C10").Value = "mROW_LOKLOD"

# C7 / C9 become explicit (blank) string cells rather than fully empty cells
$ws.Range("C7").Value = ""
$ws.Range("C9").Value = ""

# Update the saved selection to match the author's final cursor position
$ws.Range("I2:J10").Select()
